$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row values (O1 = 14, P1 = 15)
$ws.Cells.Item(1, 15).Value = 14
$ws.Cells.Item(1, 16).Value = 15

# Copy style/format from N1 (existing header style) onto O1:P1
$ws.Range("N1").Copy() | Out-Null
$ws.Range("O1:P1").PasteSpecial(-4122) | Out-Null

# Data rows 2-67, columns O (new values) and P (new values)
$ws.Cells.Item(2, 15).Value = -0.2131754727595466
$ws.Cells.Item(2, 16).Value = -0.2122923141471821
$ws.Cells.Item(3, 15).Value = 0.2801149398482115
$ws.Cells.Item(3, 16).Value = 0.2799800570166809
$ws.Cells.Item(4, 15).Value = 0.246378943620986
$ws.Cells.Item(4, 16).Value = 0.2485807460827568
$ws.Cells.Item(5, 15).Value = -0.03109275454508465
$ws.Cells.Item(5, 16).Value = -0.03262393842493257
$ws.Cells.Item(6, 15).Value = 0.2300894501463804
$ws.Cells.Item(6, 16).Value = 0.2294190496150571
$ws.Cells.Item(7, 15).Value = -0.3712679412033655
$ws.Cells.Item(7, 16).Value = -0.370134358454537
$ws.Cells.Item(8, 15).Value = -0.1424004189983353
$ws.Cells.Item(8, 16).Value = -0.1380605455470963
$ws.Cells.Item(9, 15).Value = -0.2840019017429242
$ws.Cells.Item(9, 16).Value = -0.2808933117478161
$ws.Cells.Item(10, 15).Value = 0.4300072528756786
$ws.Cells.Item(10, 16).Value = 0.4288348274975866
$ws.Cells.Item(11, 15).Value = -0.1776486411344053
$ws.Cells.Item(11, 16).Value = -0.1774130572262902
$ws.Cells.Item(12, 15).Value = -0.007330408564802555
$ws.Cells.Item(12, 16).Value = -0.01243735386936105
$ws.Cells.Item(13, 15).Value = -0.02016350281592536
$ws.Cells.Item(13, 16).Value = -0.02247199030782319
$ws.Cells.Item(14, 15).Value = 0.2055607168371144
$ws.Cells.Item(14, 16).Value = 0.1973073775944194
$ws.Cells.Item(15, 15).Value = 0.1105036265207869
$ws.Cells.Item(15, 16).Value = 0.09690696079854363
$ws.Cells.Item(16, 15).Value = 0.5597232267831413
$ws.Cells.Item(16, 16).Value = 0.5456317633757981
$ws.Cells.Item(17, 15).Value = 0.6377173923789013
$ws.Cells.Item(17, 16).Value = 0.6214248516976557
$ws.Cells.Item(18, 15).Value = -0.06822891428988065
$ws.Cells.Item(18, 16).Value = -0.07643014087704372
$ws.Cells.Item(19, 15).Value = 0.4127884575671337
$ws.Cells.Item(19, 16).Value = 0.4040932514367036
$ws.Cells.Item(20, 15).Value = 0.4606777447556887
$ws.Cells.Item(20, 16).Value = 0.4443098507845719
$ws.Cells.Item(21, 15).Value = 0.6748689376511707
$ws.Cells.Item(21, 16).Value = 0.6618296354493669
$ws.Cells.Item(22, 15).Value = 0.4442048646922097
$ws.Cells.Item(22, 16).Value = 0.43022400969447
$ws.Cells.Item(23, 15).Value = -0.00929146481804885
$ws.Cells.Item(23, 16).Value = -0.02155341910114167
$ws.Cells.Item(24, 15).Value = 2.137626993812669
$ws.Cells.Item(24, 16).Value = 1.944032336503853
$ws.Cells.Item(25, 15).Value = 0.2932176354579051
$ws.Cells.Item(25, 16).Value = 0.2911168707662191
$ws.Cells.Item(26, 15).Value = 0.1536237416503813
$ws.Cells.Item(26, 16).Value = 0.1453506594149139
$ws.Cells.Item(27, 15).Value = 0.04873815547952864
$ws.Cells.Item(27, 16).Value = 0.03971457012102407
$ws.Cells.Item(28, 15).Value = 0.7968442299292937
$ws.Cells.Item(28, 16).Value = 0.7893706035731111
$ws.Cells.Item(29, 15).Value = 1.940205889843134
$ws.Cells.Item(29, 16).Value = 1.808546135876952
$ws.Cells.Item(30, 15).Value = 0.6405966768404027
$ws.Cells.Item(30, 16).Value = 0.6336215455029707
$ws.Cells.Item(31, 15).Value = -0.4891880399900423
$ws.Cells.Item(31, 16).Value = -0.4895155768456463
$ws.Cells.Item(32, 15).Value = 0.5428480139919223
$ws.Cells.Item(32, 16).Value = 0.5363730662674459
$ws.Cells.Item(33, 15).Value = 0.7418023556976109
$ws.Cells.Item(33, 16).Value = 0.7389457296152937
$ws.Cells.Item(34, 15).Value = -0.855632259137435
$ws.Cells.Item(34, 16).Value = -0.8594499354308205
$ws.Cells.Item(35, 15).Value = 0.7743911612233859
$ws.Cells.Item(35, 16).Value = 0.77627274197615
$ws.Cells.Item(36, 15).Value = 0.712426304240861
$ws.Cells.Item(36, 16).Value = 0.7167606036379857
$ws.Cells.Item(37, 15).Value = 0.6753376977688228
$ws.Cells.Item(37, 16).Value = 0.6796718420055534
$ws.Cells.Item(38, 15).Value = 0.6244037562755153
$ws.Cells.Item(38, 16).Value = 0.6246303325296554
$ws.Cells.Item(39, 15).Value = 0.5765410465624424
$ws.Cells.Item(39, 16).Value = 0.5789048193299544
$ws.Cells.Item(40, 15).Value = 0.7396264366267948
$ws.Cells.Item(40, 16).Value = 0.741601377134533
$ws.Cells.Item(41, 15).Value = 0.5493583466030498
$ws.Cells.Item(41, 16).Value = 0.5523173013545054
$ws.Cells.Item(42, 15).Value = 0.569991015510952
$ws.Cells.Item(42, 16).Value = 0.5726449159847291
$ws.Cells.Item(43, 15).Value = 0.6501404410982843
$ws.Cells.Item(43, 16).Value = 0.6516639926053593
$ws.Cells.Item(44, 15).Value = 0.6647138508065251
$ws.Cells.Item(44, 16).Value = 0.6682344770141849
$ws.Cells.Item(45, 15).Value = 0.6213875411771471
$ws.Cells.Item(45, 16).Value = 0.6285631502390373
$ws.Cells.Item(46, 15).Value = -1.291269354579908
$ws.Cells.Item(46, 16).Value = -1.29369789774598
$ws.Cells.Item(47, 15).Value = -1.008759937911565
$ws.Cells.Item(47, 16).Value = -1.010588359964148
$ws.Cells.Item(48, 15).Value = -0.8832347647304924
$ws.Cells.Item(48, 16).Value = -0.8829727813914455
$ws.Cells.Item(49, 15).Value = -0.6446668412630039
$ws.Cells.Item(49, 16).Value = -0.6446724187429936
$ws.Cells.Item(50, 15).Value = -0.0623064789764915
$ws.Cells.Item(50, 16).Value = -0.06379448653658987
$ws.Cells.Item(51, 15).Value = -0.8708346653137617
$ws.Cells.Item(51, 16).Value = -0.8699267424614976
$ws.Cells.Item(52, 15).Value = -0.8708346653137617
$ws.Cells.Item(52, 16).Value = -0.8699267424614976
$ws.Cells.Item(53, 15).Value = -1.145456210415336
$ws.Cells.Item(53, 16).Value = -1.146220005386163
$ws.Cells.Item(54, 15).Value = -0.1771537439387354
$ws.Cells.Item(54, 16).Value = -0.1773984222435732
$ws.Cells.Item(55, 15).Value = -1.027711304077249
$ws.Cells.Item(55, 16).Value = -1.029892893037985
$ws.Cells.Item(56, 15).Value = -0.8958775851870255
$ws.Cells.Item(56, 16).Value = -0.8990681458701696
$ws.Cells.Item(57, 15).Value = -0.9039728975102416
$ws.Cells.Item(57, 16).Value = -0.9097155959343247
$ws.Cells.Item(58, 15).Value = -1.048640281008608
$ws.Cells.Item(58, 16).Value = -1.053467967632628
$ws.Cells.Item(59, 15).Value = -0.7821387775680643
$ws.Cells.Item(59, 16).Value = -0.7828368414697199
$ws.Cells.Item(60, 15).Value = -0.4040176750207024
$ws.Cells.Item(60, 16).Value = -0.405751829569123
$ws.Cells.Item(61, 15).Value = 0.3874705338414478
$ws.Cells.Item(61, 16).Value = 0.3873169753259534
$ws.Cells.Item(62, 15).Value = -1.153839959483975
$ws.Cells.Item(62, 16).Value = -1.160427858408202
$ws.Cells.Item(63, 15).Value = -0.5639347428903334
$ws.Cells.Item(63, 16).Value = -0.5601349065073576
$ws.Cells.Item(64, 15).Value = -0.839306959180219
$ws.Cells.Item(64, 16).Value = -0.839671665433906
$ws.Cells.Item(65, 15).Value = -0.0224230499572382
$ws.Cells.Item(65, 16).Value = -0.02357221422302512
$ws.Cells.Item(66, 15).Value = -0.7313674941256145
$ws.Cells.Item(66, 16).Value = -0.7379765013040298
$ws.Cells.Item(67, 15).Value = -0.6996635732273118
$ws.Cells.Item(67, 16).Value = -0.7087003905266219
